$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 161, pushing the existing rows 161:170 down to 162:171.
# Excel's Insert() carries the formatting of the row above down into the new
# row (matches the "s=2" date-format style already on column D).
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new weekly record.
$ws.Cells.Item(161, 1).Value = 8
$ws.Cells.Item(161, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(161, 3).Value = "Coquimbo"
$ws.Cells.Item(161, 4).Value = 44826
$ws.Cells.Item(161, 5).Value = 4
$ws.Cells.Item(161, 6).Value = 100112044
$ws.Cells.Item(161, 7).Value = "Perejil"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 3000
$ws.Cells.Item(161, 11).Value = 2000
$ws.Cells.Item(161, 12).Value = 2500
$ws.Cells.Item(161, 13).Value = 2250
$ws.Cells.Item(161, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(161, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(161, 16).Value = 1500
$ws.Cells.Item(161, 17).Value = 1.5
$ws.Cells.Item(161, 18).Value = "Hortaliza"
